$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Homepage-guest")

# --- Row 5 cell content updates (category/subcategory examples switched
#     from "electronics/computer/laptop" to "computers/laptop/asus") ---

$c5Text = @'
user selects a category, e.g. "computers", subcategory, e.g. "laptop" or just types the name of the product , e.g. "laptop", "asus" etc.
'@
$ws.Range("C5").Value = $c5Text

$d5Text = @'
https://www.TechShoppers.com/api/1.0/products/search?category=computers&subcategory=laptop&query=asus
'@
$ws.Range("D5").Value = $d5Text

$g5Text = @'
query - The search term e.g. "asus".
category - Search for a specific category of product (allowed values: computers, books, fashion, sports, arts & crafts)
subcategory - Search for a specific subcategory of a category of a product (allowed values: computers - laptop, desktop, pc)
'@
$ws.Range("G5").Value = $g5Text

$e5Text = @'
{
    "products": {
        "items": [
            {
                "category": "computers",
                "subcategory": "laptop",
                "brands": [
                {
                    "brand_name" : "HP",
                    "available_count": 15,
                    "platform_products": [
                    {
                        "website_name": "amazon",
                        "website_url": "https://www.amazon.com",
                        "products_info": [
                        {
                            "product_id": "2457",
                            "product_name": "HP 14 Laptop, Intel Celeron N4020",
                            "image": "https://m.media-amazon.com/images/W/MEDIAX_792452-T2/images/I/81vzIB8T1wS._AC_SL1500_.jpg",
                            "price": 60000,
                            "specs": {
                                "Processor": {
                                  "ProcessorBrand": "Intel",
                                  "ProcessorModel": "Celeron Processor N4020",
                                  "ProcessorFrequency": "1.10 GHz up to 2.80 GHz",
                                  "ProcessorCore": 2,
                                  "ProcessorThread": 2,
                                  "CPUCache": "4MB"
                                },
                                "Display": {
                                  "DisplaySize": 14,
                                  "DisplayType": "LED",
                                  "DisplayResolution": "HD (1366X768)",
                                  "TouchScreen": "No"
                                },
                                "Memory": {
                                  "RAM": "4GB(onboard)",
                                  "RAMType": "DDR4",
                                  "BusSpeed": "2666MHz"
                                }
                            },
                            "coupon" : {
                                "code":"SAVE15LAPTOP",
                                "discount_percentage": 15,
                                "discount_amount": null,
                                "description":"Save 15% on selected laptop brands",
                                "terms":"Valid on specific brands only. Excludes already discounted items",
                                "start_date":"2024-01-08",
                                "end_date":"2024-01-10"
                            },
                            "rating": 4.5
                        },
                        {
                            "product_id": "2458",
                            "product_name": "Hp 15.6\" HD Laptop Intel N200 (Pentium) Processor",
                            "image": "https://m.media-amazon.com/images/W/MEDIAX_792452-T2/images/I/51KupiNLuHL._AC_SL1280_.jpg",
                            "price": 30000,
                            "specs":{
                                "Processor": {
                                  "ProcessorBrand": "Intel",
                                  "ProcessorModel": "Core i3-1115G4",
                                  "ProcessorFrequency": "3.00 GHz up to 4.10 GHz",
                                  "ProcessorCore": 2,
                                  "ProcessorThread": 4,
                                  "CPUCache": "6 MB"
                                },
                                "Display": {
                                  "DisplaySize": "15.6 Inch",
                                  "DisplayType": "FHD LED",
                                  "DisplayResolution": "1920 x 1080",
                                  "TouchScreen": "No"
                                },
                                "Memory": {
                                  "RAM": "8GB",
                                  "RAMType": "DDR4",
                                  "BusSpeed": "3200MHz"
                                }
                            },                              
                            "coupon" : {
                                "code":"CASHBACK64",
                                "discount_percentage": null,
                                "discount_amount": null,
                                "description":"Receive BDT 6400 cashback on laptops over BDT 28000.",
                                "terms":"Cashback to be credited after purchase. Applies only to full-priced laptops.",
                                "start_date":"2024-01-09",
                                "end_date":"2024-01-11"
                            },
                            "rating": 4.4
                        }
                        ]
                    },
                    {
                        "website_name": "ebay",
                        "website_url": "https://www.ebay.com",
                        "products_info": [
                        {
                            "product_id": "2457",
                            "product_name": "HP 14 Laptop, Intel Celeron N4020",
                            "image": "https://m.media-amazon.com/images/W/MEDIAX_792452-T2/images/I/81vzIB8T1wS._AC_SL1500_.jpg",
                            "price": 65000,
                            "specs": {
                                "Processor": {
                                  "ProcessorBrand": "Intel",
                                  "ProcessorModel": "Celeron Processor N4020",
                                  "ProcessorFrequency": "1.10 GHz up to 2.80 GHz",
                                  "ProcessorCore": 2,
                                  "ProcessorThread": 2,
                                  "CPUCache": "4MB"
                                },
                                "Display": {
                                  "DisplaySize": 14,
                                  "DisplayType": "LED",
                                  "DisplayResolution": "HD (1366X768)",
                                  "TouchScreen": "No"
                                },
                                "Memory": {
                                  "RAM": "4GB(onboard)",
                                  "RAMType": "DDR4",
                                  "BusSpeed": "2666MHz"
                                }
                            },
                            "coupon" : {
                                "code":"SAVE12LAPTOP",
                                "discount_percentage": 12,
                                "discount_amount": null,
                                "description":"Save 12% on selected laptop brands",
                                "terms":"Valid on specific brands only. Excludes already discounted items",
                                "start_date":"2024-01-06",
                                "end_date":"2024-01-10"
                            },
                            "rating": 4.3
                        },
                        {
                            "product_id": "2458",
                            "product_name": "Hp 15.6\" HD Laptop Intel N200 (Pentium) Processor",
                            "image": "https://m.media-amazon.com/images/W/MEDIAX_792452-T2/images/I/51KupiNLuHL._AC_SL1280_.jpg",
                            "price": 38000,
                            "specs":{
                                "Processor": {
                                  "ProcessorBrand": "Intel",
                                  "ProcessorModel": "Core i3-1115G4",
                                  "ProcessorFrequency": "3.00 GHz up to 4.10 GHz",
                                  "ProcessorCore": 2,
                                  "ProcessorThread": 4,
                                  "CPUCache": "6 MB"
                                },
                                "Display": {
                                  "DisplaySize": "15.6 Inch",
                                  "DisplayType": "FHD LED",
                                  "DisplayResolution": "1920 x 1080",
                                  "TouchScreen": "No"
                                },
                                "Memory": {
                                  "RAM": "8GB",
                                  "RAMType": "DDR4",
                                  "BusSpeed": "3200MHz"
                                }
                            },
                            "coupon" : {
                                "code":"CASHBACK63",
                                "discount_percentage": null,
                                "discount_amount": null,
                                "description":"Receive BDT 6300 cashback on laptops over BDT 28000.",
                                "terms":"Cashback to be credited after purchase. Applies only to full-priced laptops.",
                                "start_date":"2024-01-07",
                                "end_date":"2024-01-09"
                            },
                            "rating": 4.1
                        }
                        ]
                    }
                    ]
                }, // more brands
                ]
            }
        ]
    }
}
'@
$ws.Range("E5").Value = $e5Text

# The much longer response JSON no longer fits comfortably at 14pt, so the
# author dropped E5 to the smaller (12pt) body font already used elsewhere
# in the sheet (style index 4).
$ws.Range("E5").Font.Size = 12

# Column E had to be widened considerably to accommodate the new example.
$ws.Columns("E").ColumnWidth = 144.67

# Row 9's wrapped text now auto-fits to a shorter height.
$ws.Rows("9").RowHeight = 105

# Row 5 keeps its original (already-maxed-out) custom height; re-assert it
# since the content/width edits above would otherwise let it auto-grow.
$ws.Rows("5").RowHeight = 408.6

# Selection/viewport left pointing at the edited cell.
$ws.Activate()
$ws.Range("D5").Select()
